$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: James Harden -> Keyonte George
$ws.Range("A2").Value = "Keyonte George"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Utah Jazz"

# Row 3: Bradley Beal -> Anfernee Simons
$ws.Range("A3").Value = "Anfernee Simons"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Portland Trail Blazers"

# Row 11: Moussa Diabate -> Jaxson Hayes
$ws.Range("A11").Value = "Jaxson Hayes"
$ws.Range("B11").Value = "PF,C"
$ws.Range("C11").Value = "Los Angeles Lakers"

# Row 14: Anfernee Simons -> James Harden
$ws.Range("A14").Value = "James Harden"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "LA Clippers"

# Row 15: Keyonte George -> Bradley Beal
$ws.Range("A15").Value = "Bradley Beal"
$ws.Range("B15").Value = "PG,SG,SF"
$ws.Range("C15").Value = "Phoenix Suns"
